$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for this product. It belongs right after
# the header/first couple of rows at row 9, which pushes every existing
# record at row 9 and below down by one row (old row 22 becomes row 23).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value = 44967
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(9, 6).Value = 100112001
$ws.Cells.Item(9, 7).Value = "Berenjena"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Segunda"
$ws.Cells.Item(9, 10).Value = 50
$ws.Cells.Item(9, 11).Value = 4500
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 13).Value = 4850
$ws.Cells.Item(9, 14).Value = "`$/caja 90 unidades"
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 54
$ws.Cells.Item(9, 17).Value = 90
$ws.Cells.Item(9, 18).Value = "Hortaliza"
